$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 162.5
$ws.Range("I2").Value = 162.5
$ws.Range("K2").Value = 162.5
$ws.Range("M2").Value = -49.5
$ws.Range("H12").Value = 1450
$ws.Range("I12").Value = 150
$ws.Range("J12").Value = 2750
$ws.Range("K12").Value = 150
$ws.Range("L12").Value = 2750
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = -3090
$ws.Range("H86").Value = 154324580
$ws.Range("I86").Value = 222224270
$ws.Range("K86").Value = 222224270
$ws.Range("M86").Value = -222223147
$ws.Range("H89").Value = 154324580
$ws.Range("I89").Value = 222224270
$ws.Range("K89").Value = 1111121350
$ws.Range("M89").Value = -1111115734
$ws.Range("H103").Value = 446.38
$ws.Range("I103").Value = 280.9189
$ws.Range("K103").Value = 842.7567
$ws.Range("M103").Value = -256.7567
$ws.Range("H125").Value = 62503176
$ws.Range("J125").Value = 3666.6667
$ws.Range("L125").Value = 33000.0003
$ws.Range("N125").Value = -37920.0003
$ws.Range("H141").Value = 1747.7
$ws.Range("I141").Value = 1747.7
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5243.1
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -63.10000000000036
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1472619.9
$ws.Range("I32").Value = 1507893.8
$ws.Range("K32").Value = 1507893.8
$ws.Range("M32").Value = -1507606.8
$ws.Range("H43").Value = 20000
$ws.Range("I43").Value = 20000
$ws.Range("K43").Value = 20000
$ws.Range("M43").Value = -19687
$ws.Range("H46").Value = 5333
$ws.Range("J46").Value = 5333
$ws.Range("L46").Value = 5333
$ws.Range("N46").Value = -5971
$ws.Range("H61").Value = 6069.9165
$ws.Range("I61").Value = 2771.6572
$ws.Range("K61").Value = 2771.6572
$ws.Range("M61").Value = -2559.6572
$ws.Range("H74").Value = 103441.75
$ws.Range("I74").Value = 201702.75
$ws.Range("K74").Value = 201702.75
$ws.Range("M74").Value = -200828.75
$ws.Range("H77").Value = 103441.75
$ws.Range("I77").Value = 201702.75
$ws.Range("K77").Value = 1008513.75
$ws.Range("M77").Value = -1004145.75
$ws.Range("H110").Value = 15152300
$ws.Range("I110").Value = 659.3889
$ws.Range("K110").Value = 659.3889
$ws.Range("M110").Value = 1385.6111
$ws.Range("H122").Value = 16518.555
$ws.Range("I122").Value = 27716.875
$ws.Range("K122").Value = 83150.625
$ws.Range("M122").Value = -80700.625
$ws.Range("H136").Value = 6069.9165
$ws.Range("I136").Value = 2771.6572
$ws.Range("K136").Value = 8314.971600000001
$ws.Range("M136").Value = -5764.971600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8334915.5
$ws.Range("I20").Value = 10418332
$ws.Range("J20").Value = 1247.75
$ws.Range("K20").Value = 10418332
$ws.Range("L20").Value = 1247.75
$ws.Range("M20").Value = -10418085
$ws.Range("N20").Value = -1741.75
$ws.Range("H99").Value = 2394487.8
$ws.Range("I99").Value = 1977.5714
$ws.Range("J99").Value = 5349941.5
$ws.Range("K99").Value = 1977.5714
$ws.Range("L99").Value = 5349941.5
$ws.Range("M99").Value = -479.5714
$ws.Range("N99").Value = -5352937.5
$ws.Range("H107").Value = 41672276
$ws.Range("I107").Value = 59215490
$ws.Range("J107").Value = 7150.875
$ws.Range("K107").Value = 59215490
$ws.Range("L107").Value = 7150.875
$ws.Range("M107").Value = -59213570
$ws.Range("N107").Value = -10990.875
$ws.Range("H134").Value = 7335.6333
$ws.Range("I134").Value = 2795.3333
$ws.Range("K134").Value = 8385.999899999999
$ws.Range("M134").Value = -5850.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8779.409
$ws.Range("I31").Value = 3899.0527
$ws.Range("J31").Value = 12488.48
$ws.Range("K31").Value = 3899.0527
$ws.Range("L31").Value = 12488.48
$ws.Range("M31").Value = -3604.0527
$ws.Range("N31").Value = -13078.48
$ws.Range("H34").Value = 8779.409
$ws.Range("I34").Value = 3899.0527
$ws.Range("J34").Value = 12488.48
$ws.Range("K34").Value = 3899.0527
$ws.Range("L34").Value = 12488.48
$ws.Range("M34").Value = -3697.0527
$ws.Range("N34").Value = -12892.48
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1030
$ws.Range("I3").Value = 1030
$ws.Range("K3").Value = 3090
$ws.Range("M3").Value = -2978
$ws.Range("H5").Value = 4004090.5
$ws.Range("I5").Value = 10000751
$ws.Range("J5").Value = 6316.6665
$ws.Range("K5").Value = 30002253
$ws.Range("L5").Value = 18949.9995
$ws.Range("M5").Value = -30002141
$ws.Range("N5").Value = -19173.9995
$ws.Range("H126").Value = 1900
$ws.Range("I126").Value = 1900
$ws.Range("K126").Value = 5700
$ws.Range("M126").Value = -760
$ws.Range("H135").Value = 4004090.5
$ws.Range("I135").Value = 10000751
$ws.Range("J135").Value = 6316.6665
$ws.Range("K135").Value = 90006759
$ws.Range("L135").Value = 56849.9985
$ws.Range("M135").Value = -90004224
$ws.Range("N135").Value = -61919.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H70").Value = 9309.936
$ws.Range("I70").Value = 8461.6875
$ws.Range("J70").Value = 10214.733
$ws.Range("K70").Value = 8461.6875
$ws.Range("L70").Value = 10214.733
$ws.Range("M70").Value = -8191.6875
$ws.Range("N70").Value = -10754.733
$ws.Range("H73").Value = 9309.936
$ws.Range("I73").Value = 8461.6875
$ws.Range("J73").Value = 10214.733
$ws.Range("K73").Value = 8461.6875
$ws.Range("L73").Value = 10214.733
$ws.Range("M73").Value = -7525.6875
$ws.Range("N73").Value = -12086.733
$ws.Range("H132").Value = 7450.3125
$ws.Range("I132").Value = 2764.6
$ws.Range("K132").Value = 8293.799999999999
$ws.Range("M132").Value = -5763.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8263.909
$ws.Range("J7").Value = 8778.111000000001
$ws.Range("L7").Value = 8778.111000000001
$ws.Range("N7").Value = -9002.111000000001
$ws.Range("H22").Value = 3008.5
$ws.Range("J22").Value = 3583.5715
$ws.Range("L22").Value = 3583.5715
$ws.Range("N22").Value = -4173.5715
$ws.Range("H27").Value = 3008.5
$ws.Range("J27").Value = 3583.5715
$ws.Range("L27").Value = 3583.5715
$ws.Range("N27").Value = -3797.5715
$ws.Range("H40").Value = 5777.6763
$ws.Range("I40").Value = 4404.6313
$ws.Range("K40").Value = 4404.6313
$ws.Range("M40").Value = -4268.6313
$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 4000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -4376
$ws.Range("H55").Value = 58824120
$ws.Range("I55").Value = 250000100
$ws.Range("J55").Value = 743.1539
$ws.Range("K55").Value = 250000100
$ws.Range("L55").Value = 743.1539
$ws.Range("M55").Value = -249999927
$ws.Range("N55").Value = -1089.1539
$ws.Range("H61").Value = 4372.852
$ws.Range("I61").Value = 1269.7858
$ws.Range("K61").Value = 1269.7858
$ws.Range("M61").Value = -1067.7858
$ws.Range("H68").Value = 6578.5713
$ws.Range("I68").Value = 6000
$ws.Range("J68").Value = 6810
$ws.Range("K68").Value = 6000
$ws.Range("L68").Value = 6810
$ws.Range("M68").Value = -5251
$ws.Range("N68").Value = -8308
$ws.Range("H71").Value = 6578.5713
$ws.Range("I71").Value = 6000
$ws.Range("J71").Value = 6810
$ws.Range("K71").Value = 30000
$ws.Range("L71").Value = 34050
$ws.Range("M71").Value = -26256
$ws.Range("N71").Value = -41538
$ws.Range("H113").Value = 4372.852
$ws.Range("I113").Value = 1269.7858
$ws.Range("K113").Value = 1269.7858
$ws.Range("M113").Value = 900.2141999999999
$ws.Range("H126").Value = 8263.909
$ws.Range("J126").Value = 8778.111000000001
$ws.Range("L126").Value = 26334.333
$ws.Range("N126").Value = -31274.333
$ws.Range("H132").Value = 7941876.5
$ws.Range("I132").Value = 13515962
$ws.Range("K132").Value = 40547886
$ws.Range("M132").Value = -40545356

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H107").Value = 15152389
$ws.Range("I107").Value = 704
$ws.Range("J107").Value = 47620284
$ws.Range("K107").Value = 2112
$ws.Range("L107").Value = 142860852
$ws.Range("M107").Value = -192
$ws.Range("N107").Value = -142864692
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H126").Value = 3175.1538
$ws.Range("J126").Value = 7499.75
$ws.Range("L126").Value = 22499.25
$ws.Range("N126").Value = -27439.25
$ws.Range("H132").Value = 15629761
$ws.Range("I132").Value = 22732656
$ws.Range("K132").Value = 68197968
$ws.Range("M132").Value = -68195438
